$d = $word.ActiveDocument

# 1. Replace the wording of the "December campaigns" bullet with the new sentence.
$d.Content.Find.Execute(
    "Campaigns that were created in December were the least likely to be successful",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "More campaigns that were created in December failed than were successful", 2)

# 2. Move the "_GoBack" bookmark from the title paragraph to the start of the
#    "The categories with the highest rate of success" bullet paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$target = $d.Paragraphs(8).Range
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)

Write-Host "done"
